# edit.ps1 — apply the "Code metric.docx" revision described by the commit
# "Tilføjet lidt mere til Code metric.docx".
#
# The heading's visible text is unchanged ("Code metric" stays "Code metric" —
# the source diff only wraps "metric" in spell-check proofing marks, which
# carry no visible content and aren't something the Word object model exposes
# for scripting, so the heading paragraph is left untouched).
#
# The body paragraph gains extra wording / punctuation in several spots. We
# apply a sequence of small, precisely-anchored Find & Replace operations —
# each `Find` string unique in the document — working from the end of the
# paragraph back towards the start, so every edit only folds together the
# run(s) it actually touches instead of needlessly disturbing text that
# comes before it.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute(
        $find,    # FindText
        $true,    # MatchCase
        $false,   # MatchWholeWord
        $false,   # MatchWildcards
        $false,   # MatchSoundsLike
        $false,   # MatchAllWordForms
        $true,    # Forward
        1,        # Wrap (wdFindContinue)
        $false,   # Format
        $replace, # ReplaceWith
        2         # Replace (wdReplaceAll)
    ) | Out-Null
}

# 10) "Det har medført at projektets" -> "Det har medført, at projektets"
Replace-Text `
    "Det har medført at projektets" `
    "Det har medført, at projektets"

# 9) "cyklomatisk" -> "cyklomatiske" (3rd occurrence)
Replace-Text `
    "holde den cyklomatisk kompleksitet nede, ved" `
    "holde den cyklomatiske kompleksitet nede, ved"

# 8) "cyklomatisk" -> "cyklomatiske" (2nd occurrence)
Replace-Text `
    "den cyklomatisk kompleksitet nede. Der" `
    "den cyklomatiske kompleksitet nede. Der"

# 7) "vedligeholdsvenlig, som mulig" -> "vedligeholdsvenlig som mulig" (drop comma)
Replace-Text `
    "vedligeholdsvenlig, som mulig" `
    "vedligeholdsvenlig som mulig"

# 6) new passage inserted right after "eller refaktoreres." and before
#    "Der er gennem projektet forsøgt ..."
Replace-Text `
    "eller refaktoreres. Der er gennem projektet forsøgt" `
    "eller refaktoreres. Det gør det også muligt at identificere de funktioner, der har en høj potentiel risiko. Derudfra kan der laves en vurdering på, om arbejdet med at ændre funktionen så risikoen for den bliver mindre, giver en stor nok fordel i forhold til den tid, der skal bruges på det. Der kan også besluttes, at funktionen bare skal testes godt igennem på grund af, at den har en øget risiko. Der er gennem projektet forsøgt"

# 5) "et værktøj til" -> "et værktøj, til"
Replace-Text `
    "et værktøj til" `
    "et værktøj, til"

# 4) "vedligeholdelsesvenligt og" -> "vedligeholdelsesvenlig, og"
Replace-Text `
    "vedligeholdelsesvenligt og" `
    "vedligeholdelsesvenlig, og"

# 3) "Dette er gjort for at koden" -> "Dette er gjort for, at koden"
Replace-Text `
    "gjort for at koden" `
    "gjort for, at koden"

# 2) "code metric." -> "code metrics."
Replace-Text `
    "code metric. Dette" `
    "code metrics. Dette"

# 1) "cyklomatisk" -> "cyklomatiske" (1st occurrence)
Replace-Text `
    "cyklomatisk kompleksitet for de forskellige" `
    "cyklomatiske kompleksitet for de forskellige"
